$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Precondition rows: fill in "Есть ..." results (column B) for rows 7-10 ---
$ws.Range("B7").Value = "Есть шариковая механическаяручка."
$ws.Range("B8").Value = "Есть бумага."
$ws.Range("B9").Value = "Есть письменный стол."
$ws.Range("B10").Value = "Есть  стул."

# --- Steps block (rows 15-18): re-flow / correct the step texts ---
$ws.Range("A15").Value = "1.Взять ручку и бумагу, положить на стол."
$ws.Range("B15").Value = "Бумага и ручка лежит на столе"

$ws.Range("A16").Value = "2.Сесть на стул за стол."
$ws.Range("B16").Value = "Сидим на стуле за столом, на столе ручка и бумага"

$ws.Range("A17").Value = "3. Включить ручку."
$ws.Range("B17").Value = "Ручка в рабочем положении"

# Row 18 used to be blank - now it becomes part of the steps table with its own
# merged B:C cell and bottom "Пройден" status, matching rows 15-17.
$ws.Range("B18:C18").Merge()
$ws.Range("A18").Value = "4. Написать на бумаге 50 000 слов."
$ws.Range("B18").Value = "Написали 50 010 слов"
$ws.Range("D18").Value = "Пройден"

# Normalize borders for the whole steps block (rows 15-18) so the newly-added
# row 18 matches the thin full box border used by the rest of the table.
$ws.Range("A15:D18").Borders.LineStyle = 1
$ws.Range("A15:D18").Borders.Weight = 2

# --- Postcondition rows (24): fill in the "switch off the pen" step ---
$ws.Range("A24").Value = "Выключить ручку"
$ws.Range("B24").Value = "Ручка в нерабочем положении"
$ws.Range("D24").Value = "Пройден"

# --- Selection moves from the old B17:C17 to the newly populated A18 ---
$ws.Range("A18").Select()
